$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new time-tracking entries (date, hours, description) to the
# "Jiska" section (columns M/N/O), rows 20-21.
$ws.Range("M20").Value = 45356
$ws.Range("N20").Value = 3
$ws.Range("O20").Value = "Ohjelmiston kehitys, taulukko, printtaus sekä esittäminen"

$ws.Range("M21").Value = 45357
$ws.Range("N21").Value = 2
$ws.Range("O21").Value = "Ohjelmiston kehity, merge ongelman setviminen"

# Adjust the view: zoom out to 70% and move the selection to O24.
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("O24").Select()
